$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Various combinations of closing the window and re-opening"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Interaction with LCS manager and EAP"

$ws.Range("C8").Select()
